# Update the "报名人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 265
    3  = 444
    4  = 13637
    7  = 43
    9  = 173
    11 = 484
    13 = 79
    17 = 440
    18 = 5650
    20 = 75
    21 = 968
    22 = 52
    24 = 137
    25 = 186
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
